# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing a Text number format so
# numeric-looking strings (e.g. "217.50", "1.003") are not silently coerced
# into floating point numbers by the COM Value setter (which would drop
# trailing zeros / change the stored cell type). Style is reset back to
# "Normal" afterwards so the cell keeps its original (unstyled) appearance.
function Set-TextCell {
    param($ws, $addr, $value)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "26.101.85"
Set-TextCell $ws "E2" "  -0.28%  "

Set-TextCell $ws "D3" "1.652.59"
Set-TextCell $ws "E3" "  -0.72%  "

Set-TextCell $ws "D5" "217.50"

Set-TextCell $ws "D6" "0.5280"
Set-TextCell $ws "E6" "  +0.47%  "

Set-TextCell $ws "D7" "1.003"
Set-TextCell $ws "E7" "  -0.02%  "

Set-TextCell $ws "D8" "0.2605"
Set-TextCell $ws "E8" "  -1.52%  "

Set-TextCell $ws "E9" "  +0.52%  "

Set-TextCell $ws "D10" "20.33"
Set-TextCell $ws "E10" "  -2.00%  "

Set-TextCell $ws "D11" "0.07792"
Set-TextCell $ws "E11" "  +0.58%  "

Set-TextCell $ws "D12" "4.521"
Set-TextCell $ws "E12" "  +1.22%  "

Set-TextCell $ws "D13" "1.651.03"
Set-TextCell $ws "E13" "  +0.47%  "

Set-TextCell $ws "D14" "1.880.49"
Set-TextCell $ws "E14" "  -0.62%  "

Set-TextCell $ws "D15" "0.5472"
Set-TextCell $ws "E15" "  -0.03%  "

Set-TextCell $ws "D16" "0.0₅8190"
Set-TextCell $ws "E16" "  +0.83%  "

Set-TextCell $ws "D17" "65.37"
Set-TextCell $ws "E17" "  +0.71%  "

Set-TextCell $ws "D18" "26.101.00"
Set-TextCell $ws "E18" "  -0.30%  "

Set-TextCell $ws "D19" "1.003"
Set-TextCell $ws "E19" "  +0.02%  "

Set-TextCell $ws "D20" "4.581"
Set-TextCell $ws "E20" "  -0.29%  "

Set-TextCell $ws "D21" "190.65"
Set-TextCell $ws "E21" "  -0.69%  "

Set-TextCell $ws "E22" "  +0.31%  "

Set-TextCell $ws "D23" "6.011"
Set-TextCell $ws "E23" "  -0.01%  "

Set-TextCell $ws "D24" "1.004"
Set-TextCell $ws "E24" "  -0.03%  "

Set-TextCell $ws "D25" "144.83"

Set-TextCell $ws "D26" "0.1226"
Set-TextCell $ws "E26" "  -1.06%  "

Set-TextCell $ws "D27" "7.213"
Set-TextCell $ws "E27" "  -0.67%  "

Set-TextCell $ws "E28" "  -1.52%  "

Set-TextCell $ws "D29" "1.453"
Set-TextCell $ws "E29" "  +3.74%  "

Set-TextCell $ws "D30" "0.05771"
Set-TextCell $ws "E30" "  -3.48%  "

Set-TextCell $ws "D31" "1.271"
Set-TextCell $ws "E31" "  -0.72%  "

Set-TextCell $ws "D32" "3.545"
Set-TextCell $ws "E32" "  +0.29%  "

Set-TextCell $ws "D33" "3.260"
Set-TextCell $ws "E33" "  -0.28%  "

Set-TextCell $ws "D34" "1.596"
Set-TextCell $ws "E34" "  +1.11%  "

Set-TextCell $ws "D35" "2.800"
Set-TextCell $ws "E35" "  +1.13%  "

Set-TextCell $ws "D36" "2.415"
Set-TextCell $ws "E36" "  -0.03%  "

Set-TextCell $ws "D37" "0.9453"
Set-TextCell $ws "E37" "  -1.69%  "

Set-TextCell $ws "D38" "0.5746"
Set-TextCell $ws "E38" "  +1.35%  "

Set-TextCell $ws "D39" "0.01611"
Set-TextCell $ws "E39" "  +0.77%  "

Set-TextCell $ws "D40" "0.8519"
Set-TextCell $ws "E40" "  +0.01%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D41" "104.32"
Set-TextCell $ws "E41" "  +3.19%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D42" "104.32"
Set-TextCell $ws "E42" "  +0.05%  "

Set-TextCell $ws "D43" "5.716"
Set-TextCell $ws "E43" "  -3.77%  "

Set-TextCell $ws "D44" "1.030.89"
Set-TextCell $ws "E44" "  +2.47%  "

Set-TextCell $ws "D45" "1.794.18"
Set-TextCell $ws "E45" "  -0.70%  "

Set-TextCell $ws "D46" "56.87"
Set-TextCell $ws "E46" "  +0.06%  "

Set-TextCell $ws "D47" "1.001"
Set-TextCell $ws "E47" "  +0.35%  "

Set-TextCell $ws "D48" "0.4332"
Set-TextCell $ws "E48" "  +0.38%  "

Set-TextCell $ws "D49" "7.855"
Set-TextCell $ws "E49" "  -1.80%  "

Set-TextCell $ws "E50" "  -0.22%  "

Set-TextCell $ws "D51" "1.446"
Set-TextCell $ws "E51" "  -1.57%  "
